$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix objective-function bug ---
# Month 11 (row 24) inputs were missing a planned/repair event; correct the
# recorded n0/n1/Repair counts so the downstream Leftover hours, Charter
# cost, Breakage cost etc. (all formulas) recalculate correctly.
$ws.Range("C24").Value = 1
$ws.Range("D24").Value = 1
$ws.Range("F24").Value = 3

# --- Remove the scratch "FIND OUT DIFFERENCES" check columns/cells ---
# These were a temporary cross-check (N/O helper columns + the C10/D10
# comparison formula) that is no longer needed now that the bug is fixed.
$ws.Range("C10:D10").Clear()
$ws.Range("N12:O26").Clear()

# --- Rename the "Broken at end" column header ---
$ws.Range("K12").Value = "Fails at end"

# --- Column widths (best-effort cosmetic match) ---
$ws.Columns(11).ColumnWidth = 10.17
$ws.Columns(14).ColumnWidth = 12.6

# --- Restore selection to the data-entry block ---
$ws.Range("C13:F24").Select() | Out-Null
